$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the newly-added availability data for Thursday (row 4)
# (order matches the shared-string table build order from the original edit)
$ws.Range("B4").Value = "9am-3pm"
$ws.Range("C4").Value = "11am-3pm"
$ws.Range("D4").Value = "9am-3pm"
$ws.Range("E4").Value = "11am-3pm"
$ws.Range("H4").Value = "10am-6pm"
$ws.Range("G4").Value = "10am-MN"
$ws.Range("F4").Value = "9am-MN"

# Match the column C width change recorded in the saved file
$ws.Columns.Item(3).ColumnWidth = 8.3

# Move the active selection to match the saved view state
$ws.Range("F17").Select()
